# Apply "valueset and codesystem changes" edit:
#  - Metadata sheet, URL row value changes to the HL7 immunization-status URL
#  - Metadata sheet, Date row value is bumped to the newer timestamp
#  - The "Include #0" sheet's "System URI" value (already the HL7 URL) is left as-is;
#    it ends up sharing the same string table entry as the updated URL above.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")

$metadata.Range("B2").Value = "http://hl7.org/fhir/ValueSet/immunization-status"
$metadata.Range("B8").Value = "2025-06-25T06:29:04+01:00"
